# "updates to data dictionary"
#
# Mark the "_id" and "Division" rows for deletion (col C = "X"), mark the
# "Neighbourhood" row for deletion too (col C = "x", matching the existing
# lower-case marker already used elsewhere in the sheet), and clear the
# green "needs review" highlight from the Description cells for the
# "Division" and "Neighbourhood" rows now that they've been filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  ("_id")          -> flag column C for deletion
$ws.Range("C2").Value = "X"

# Row 10 ("Division")     -> flag column C for deletion, un-highlight B10
$ws.Range("C10").Value = "X"
$ws.Range("B10").Interior.ColorIndex = -4142

# Row 12 ("Neighbourhood")-> flag column C for deletion, un-highlight B12
$ws.Range("C12").Value = "x"
$ws.Range("B12").Interior.ColorIndex = -4142

# Leave the selection where the edits ended
$ws.Range("C9").Select()
